$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the HTTP row (row 7) results
$ws.Range("B7").Value = 57.23
$ws.Range("C7").Value = 10.1
$ws.Range("D7").Value = 2791.76
$ws.Range("E7").Value = 380.36
$ws.Range("F7").Value = 56156.39
$ws.Range("G7").Value = 9694.0499999999993
$ws.Range("H7").Value = 67710.94
$ws.Range("I7").Value = 7408.72
$ws.Range("K7").Value = 1.0198199999999999
$ws.Range("L7").Value = 1.0002
$ws.Range("M7").Value = 1.0000199999999999

# Update the active selection on the sheet
$ws.Range("M8").Select()
